# Workbook: mapping_docs/cars-com_card mapping.xlsx
# Commit: "updated procedures - full reload data from Landing to Stg1,
#          usp_write_event_log for logging procedures that were run"
#
# Logical edit performed on the "Map lnd_cars-com_card to Stg" sheet:
#   1. Remove the obsolete "dl_loaded_date" mapping row (was row 34) - this
#      column/row is no longer populated because Landing -> Stg1 is now a
#      full reload rather than an incremental one, so the old
#      CURRENT_TIMESTAMP BETWEEN date-range example row goes away too.
#   2. Rename the "stg1_loaded_date" field to "modified_date" (now set by
#      usp_write_event_log), in what becomes row 34 after the deletion.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Map lnd_cars-com_card to Stg")
$ws.Activate()

# 1) Delete the whole "dl_loaded_date" row - shifts the two rows below it
#    (stg1_loaded_date / row_hash) up by one.
$ws.Rows.Item(34).Delete()

# 2) Rename "stg1_loaded_date" -> "modified_date" in the column that now
#    sits on row 34.
$ws.Range("D34").Value = "modified_date"

# Leave the selection where the author left it when they saved the file.
$ws.Range("H36").Select()
